# Apply attendance-count updates to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H ("Absent") updates: 0 -> 1
$ws.Range("H3").Value  = 1
$ws.Range("H4").Value  = 1
$ws.Range("H5").Value  = 1
$ws.Range("H6").Value  = 1
$ws.Range("H7").Value  = 1
$ws.Range("H8").Value  = 1
$ws.Range("H9").Value  = 1
$ws.Range("H11").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("H18").Value = 1

# Row 10 updates: D (Total Attendance Count), E (Real), F (Duplicate)
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1

# Row 12 updates: D, E
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

# Row 13 updates: D, E, F
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 1

# Row 14 updates: D, E, F
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 1

# Row 15 updates: D, E, F
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 1

# Row 17 updates: D, E, F
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 1
